$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "26.290.80"
$ws.Range("E2").Value = "  +0.99%  "
Set-TextValue $ws.Range("D3") "1.679.23"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  +0.37%  "
Set-TextValue $ws.Range("D5") "217.98"
$ws.Range("E5").Value = "  +0.77%  "
Set-TextValue $ws.Range("D6") "0.5262"
Set-TextValue $ws.Range("D7") "1.008"
$ws.Range("E7").Value = "  +0.34%  "
Set-TextValue $ws.Range("D8") "0.2686"
$ws.Range("E8").Value = "  +2.22%  "
Set-TextValue $ws.Range("D9") "0.06466"
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("E10").Value = "  +1.07%  "
Set-TextValue $ws.Range("D11") "0.07505"
$ws.Range("E11").Value = "  +1.17%  "
Set-TextValue $ws.Range("D12") "1.699.14"
$ws.Range("E12").Value = "  +1.88%  "
Set-TextValue $ws.Range("D13") "4.513"
$ws.Range("E13").Value = "  +0.41%  "
Set-TextValue $ws.Range("D14") "0.5777"
$ws.Range("E14").Value = "  -0.41%  "
Set-TextValue $ws.Range("D15") "0.000008505"
$ws.Range("E15").Value = "  -0.40%  "
Set-TextValue $ws.Range("D16") "64.71"
$ws.Range("E16").Value = "  +0.81%  "
Set-TextValue $ws.Range("D17") "26.346.73"
$ws.Range("E17").Value = "  +1.01%  "
Set-TextValue $ws.Range("D18") "4.916"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("E19").Value = "  +0.39%  "
Set-TextValue $ws.Range("D20") "10.87"
$ws.Range("E20").Value = "  +1.55%  "
Set-TextValue $ws.Range("D21") "189.58"
$ws.Range("E21").Value = "  +0.49%  "
Set-TextValue $ws.Range("D22") "6.188"
$ws.Range("E23").Value = "  +0.31%  "
Set-TextValue $ws.Range("D24") "144.93"
$ws.Range("E24").Value = "  -0.58%  "
Set-TextValue $ws.Range("D25") "7.778"
$ws.Range("E25").Value = "  +2.42%  "
$ws.Range("E26").Value = "  +5.69%  "
$ws.Range("E27").Value = "  +1.09%  "
Set-TextValue $ws.Range("D28") "0.06445"
$ws.Range("E28").Value = "  -1.74%  "
Set-TextValue $ws.Range("D29") "1.363"
$ws.Range("E29").Value = "  +4.13%  "
$ws.Range("E30").Value = "  +0.86%  "
Set-TextValue $ws.Range("D31") "3.585"
$ws.Range("E31").Value = "  +1.79%  "
Set-TextValue $ws.Range("D32") "3.585"
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("E33").Value = "  +2.08%  "
Set-TextValue $ws.Range("D34") "1.028"
$ws.Range("E34").Value = "  +1.04%  "
Set-TextValue $ws.Range("D35") "0.6201"
$ws.Range("E35").Value = "  +2.49%  "
Set-TextValue $ws.Range("D37") "2.740"
$ws.Range("E37").Value = "  +2.12%  "
Set-TextValue $ws.Range("D38") "6.276"
$ws.Range("E38").Value = "  +1.27%  "
Set-TextValue $ws.Range("D39") "1.116.85"
$ws.Range("E39").Value = "  +3.92%  "
$ws.Range("E40").Value = "  +0.65%  "
Set-TextValue $ws.Range("D41") "0.8730"
$ws.Range("E41").Value = "  +1.49%  "
Set-TextValue $ws.Range("D42") "1.015"
$ws.Range("E42").Value = "  +0.65%  "
Set-TextValue $ws.Range("D43") "100.52"
$ws.Range("E43").Value = "  +0.04%  "
Set-TextValue $ws.Range("D44") "1.828.66"
$ws.Range("E44").Value = "  +0.90%  "
Set-TextValue $ws.Range("D45") "0.00000000109"
$ws.Range("E45").Value = "  -3.64%  "
$ws.Range("E47").Value = "  +2.12%  "
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("E49").Value = "  +1.18%  "
Set-TextValue $ws.Range("D50") "0.4293"
$ws.Range("E50").Value = "  +0.07%  "
Set-TextValue $ws.Range("D51") "6.057"
$ws.Range("E51").Value = "  +1.92%  "
